# Germany Regionalliga West - base update (30-03-2024)
# Applies:
#  1) A set of row-data rotations/swaps among existing match rows
#     (columns B:AC move between rows; column A, the sequential
#     index, stays put on each row).
#  2) Full replacement of rows 235 and 236 with newly scraped
#     fixture data.
#  3) Deletion of the trailing placeholder rows 237-240.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row-data permutations -------------------------------------------
# Each inner array is a cycle of row numbers: the data (columns B:AC)
# that currently lives in row group[i+1] ends up in row group[i], the
# last element's data wraps around to the first.
$groups = @(
    @(19,20),
    @(40,42,41),
    @(82,83,84,85,86),
    @(108,109),
    @(125,126),
    @(128,129),
    @(153,154),
    @(161,162),
    @(166,168),
    @(205,209),
    @(216,221,220,219,218,217)
)

foreach ($group in $groups) {
    $n = $group.Length

    # Snapshot the current (pre-edit) values for every row in the cycle
    # before writing anything, since several rows feed into each other.
    $snapshots = @()
    foreach ($r in $group) {
        $snapshots += ,($ws.Range("B$r`:AC$r").Value2)
    }

    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcValues = $snapshots[($i + 1) % $n]
        $ws.Range("B$destRow`:AC$destRow").Value2 = $srcValues
    }
}

# --- 2) Replace rows 235 & 236 with new fixture data ---------------------
$ws.Range("B235").Value2 = 6886997
$ws.Range("E235").Value2 = 45382.375
$ws.Range("F235").Value2 = "SC Paderborn 07 II"
$ws.Range("G235").Value2 = "SV Rodinghausen"
$ws.Range("K235").Value2 = 2.7
$ws.Range("L235").Value2 = 3.6
$ws.Range("M235").Value2 = 2.15
$ws.Range("N235").Value2 = 3.4
$ws.Range("O235").Value2 = 3.6
$ws.Range("P235").Value2 = 1.85
$ws.Range("Q235").Value2 = 0.5
$ws.Range("R235").Value2 = 1.875
$ws.Range("S235").Value2 = 1.975
$ws.Range("T235").Value2 = 2.75
$ws.Range("U235").Value2 = 1.975
$ws.Range("V235").Value2 = 1.875

$ws.Range("B236").Value2 = 6886996
$ws.Range("E236").Value2 = 45382.375
$ws.Range("F236").Value2 = "FC WegbergBeeck"
$ws.Range("G236").Value2 = "Cologne II"
$ws.Range("K236").Value2 = 3
$ws.Range("L236").Value2 = 3.6
$ws.Range("M236").Value2 = 2
$ws.Range("N236").Value2 = 3.6
$ws.Range("O236").Value2 = 3.8
$ws.Range("P236").Value2 = 1.75
$ws.Range("Q236").Value2 = 0.75
$ws.Range("R236").Value2 = 1.825
$ws.Range("S236").Value2 = 2.025
$ws.Range("T236").Value2 = 3
$ws.Range("U236").Value2 = 1.95
$ws.Range("V236").Value2 = 1.9

# --- 3) Drop the trailing placeholder rows 237-240 ------------------------
$ws.Range("A237:AC240").EntireRow.Delete() | Out-Null
